$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 863
$ws1.Range("F3").Value = 13808
$ws1.Range("F4").Value = 13598
$ws1.Range("F8").Value = 599
$ws1.Range("F11").Value = 56
$ws1.Range("F12").Value = 765
$ws1.Range("F14").Value = 108
$ws1.Range("F17").Value = 125
$ws1.Range("F19").Value = 531
$ws1.Range("F21").Value = 408
$ws1.Range("F22").Value = 325
$ws1.Range("F23").Value = 266
$ws1.Range("F25").Value = 90
$ws1.Range("F26").Value = 5

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 166
$ws2.Range("F7").Value = 1523
$ws2.Range("F11").Value = 68

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 111

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 863
$ws4.Range("F4").Value = 13808
$ws4.Range("F5").Value = 13598
$ws4.Range("F9").Value = 599
$ws4.Range("F12").Value = 56
$ws4.Range("F13").Value = 765
$ws4.Range("F17").Value = 108
$ws4.Range("F20").Value = 125
$ws4.Range("F24").Value = 111
$ws4.Range("F25").Value = 111
$ws4.Range("F26").Value = 531
$ws4.Range("F28").Value = 408
$ws4.Range("F29").Value = 325
$ws4.Range("F30").Value = 266
$ws4.Range("F32").Value = 166
$ws4.Range("F33").Value = 1523
$ws4.Range("F37").Value = 90
$ws4.Range("F38").Value = 68
$ws4.Range("F40").Value = 5
